$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 145, shifting existing rows 145-195 down to 146-196
$ws.Rows.Item(145).EntireRow.Insert()

# Populate the newly inserted row 145 with the new record
$ws.Cells.Item(145, 1).Value = 8
$ws.Cells.Item(145, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(145, 3).Value = "Coquimbo"
$ws.Cells.Item(145, 4).Value = (Get-Date -Year 2021 -Month 9 -Day 29 -Hour 0 -Minute 0 -Second 0).Date
$ws.Cells.Item(145, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(145, 5).Value = 4
$ws.Cells.Item(145, 6).Value = 100114013
$ws.Cells.Item(145, 7).Value = "Zanahoria"
$ws.Cells.Item(145, 8).Value = "Sin especificar"
$ws.Cells.Item(145, 9).Value = "Primera"
$ws.Cells.Item(145, 10).Value = 800
$ws.Cells.Item(145, 11).Value = 6000
$ws.Cells.Item(145, 12).Value = 7000
$ws.Cells.Item(145, 13).Value = 6500
$ws.Cells.Item(145, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(145, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(145, 16).Value = 325
$ws.Cells.Item(145, 17).Value = 20
$ws.Cells.Item(145, 18).Value = "Hortaliza"
